# Development_Roadmap.xlsx -- "Documentation of V0.7.0 beta finished"
#
# 1. Insert a new row (job-control bug fix) after the "Algorithm doc update" header row.
# 2. Insert a second new row (manual frame exclusion feature) further down.
# 3. Mark three "Must have" documentation tasks as "done".
# 4. Append two blank rows at the bottom of the table (handled automatically by the
#    row-insert shifts below).
# 5. Update the active selection shown when the sheet is opened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---------------------------------------------------------------------------
# Step 1: insert a new row at position 4 describing the job-control bug fix.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).EntireRow.Insert()

$ws.Cells.Item(4,1).Value = 'Bug if the user wants to go back to a previous step after a job is aborted'
$ws.Cells.Item(4,2).Value = 'If a job is aborted because of a runtime error, execution continues with the next job, or the program is set to idle mode (if there are no more jobs). If then the user presses "go back to" he is presented with the full set of job steps to choose from, even if they have not been executed for the aborted job. If he then chooses a step which was not executed, the program crashes without further notice.'
$ws.Cells.Item(4,3).Value = 'The job control in the main GUI thread must be corrected.'
$ws.Cells.Item(4,5).Value = 'Rolf'
$ws.Cells.Item(4,6).Value = 'Bug fix'
$ws.Cells.Item(4,7).Value = '0.7.0'
$ws.Cells.Item(4,8).Value = 'open'
$ws.Rows.Item(4).RowHeight = 90

# ---------------------------------------------------------------------------
# Step 2: mark the three "Must have" documentation rows as "done" (these rows
# have shifted down by one because of the insert above: old rows 7,8,9 are now
# rows 8,9,10).
# ---------------------------------------------------------------------------
$ws.Cells.Item(8,8).Value  = 'done'
$ws.Cells.Item(9,8).Value  = 'done'
$ws.Cells.Item(10,8).Value = 'done'

# ---------------------------------------------------------------------------
# Step 3: insert a new row at position 13 describing the manual frame
# exclusion feature (old row 12 "Add frame stabilization None" now sits at
# row 13, so the new row is inserted right before it).
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).EntireRow.Insert()

$ws.Cells.Item(13,1).Value = 'Add the option to manually exclude frames from the input stack'
$ws.Cells.Item(13,2).Value = 'This should be added to the phase where the stack size is set. When the user scrolls through the video (using the slider or the number boxes), an additional checkbox can be set / unset to include / exclude a given frame from the input stack.'
$ws.Cells.Item(13,3).Value = 'The implementation should be done in the module "frames" via an index translation table. This way (by going back to this step) frames can be restored by simply resetting the translation table.'
$ws.Cells.Item(13,5).Value = 'Rolf'
$ws.Cells.Item(13,6).Value = 'Must have'
$ws.Cells.Item(13,7).Value = '> 0.7.0'
$ws.Cells.Item(13,8).Value = 'open'
$ws.Rows.Item(13).RowHeight = 75

# ---------------------------------------------------------------------------
# Step 4: update the selection stored with the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("A1:C1").Select()

Write-Output "Development roadmap updated for V0.7.0 beta documentation."
